# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the "展览" (Exhibition), "演出" (Performance), and "全部类型" (All types)
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 220
$ws1.Range("F3").Value = 526
$ws1.Range("F7").Value = 65
$ws1.Range("F8").Value = 7104
$ws1.Range("F10").Value = 409
$ws1.Range("F11").Value = 3568
$ws1.Range("F12").Value = 301
$ws1.Range("F13").Value = 521
$ws1.Range("F14").Value = 264
$ws1.Range("F15").Value = 591
$ws1.Range("F16").Value = 70

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = 178

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 178
$ws4.Range("F4").Value = 220
$ws4.Range("F5").Value = 526
$ws4.Range("F9").Value = 65
$ws4.Range("F11").Value = 7104
$ws4.Range("F14").Value = 409
$ws4.Range("F15").Value = 3568
$ws4.Range("F16").Value = 301
$ws4.Range("F17").Value = 521
$ws4.Range("F18").Value = 264
$ws4.Range("F19").Value = 591
$ws4.Range("F20").Value = 70
